$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.06522651944612287
$ws.Range("C2").Value = 1.528336551320912
$ws.Range("D2").Value = 13.19799231504191
$ws.Range("E2").Value = 3.632904115861292
$ws.Range("F2").Value = 3.667756611022197

$ws.Range("B3").Value = 0.4071403425012899
$ws.Range("C3").Value = 1.360982216433359
$ws.Range("D3").Value = 11.3091244196728
$ws.Range("E3").Value = 3.362904164509122
$ws.Range("F3").Value = 3.371383749205116

$ws.Range("B4").Value = 0.5746585884689557
$ws.Range("C4").Value = 1.326100425856204
$ws.Range("D4").Value = 8.342275789519231
$ws.Range("E4").Value = 2.888299809493334
$ws.Range("F4").Value = 2.859292636981395

$ws.Range("B5").Value = 0.5077119487392797
$ws.Range("C5").Value = 1.431591456579723
$ws.Range("D5").Value = 7.506921903539223
$ws.Range("E5").Value = 2.739876256975709
$ws.Range("F5").Value = 2.725872198519724
$ws.Range("G5").Value = 41

$ws.Range("B6").Value = 0.5683032239805089
$ws.Range("C6").Value = 1.971558803322826
$ws.Range("D6").Value = 12.24721082536355
$ws.Range("E6").Value = 3.499601523797181
$ws.Range("F6").Value = 3.510230336223366
$ws.Range("G6").Value = 31

$ws.Range("B7").Value = 0.5396674267054108
$ws.Range("C7").Value = 2.030587006868052
$ws.Range("D7").Value = 12.72930806051201
$ws.Range("E7").Value = 3.56781558667373
$ws.Range("F7").Value = 3.589189854583631
$ws.Range("G7").Value = 29

$ws.Range("B8").Value = 0.4461851316959971
$ws.Range("C8").Value = 2.04993939170431
$ws.Range("D8").Value = 12.97625138835586
$ws.Range("E8").Value = 3.602256430122079
$ws.Range("F8").Value = 3.642608933213258
$ws.Range("G8").Value = 27

$ws.Range("B9").Value = 0.2350498390188732
$ws.Range("C9").Value = 2.433167982692137
$ws.Range("D9").Value = 17.26964117102916
$ws.Range("E9").Value = 4.155675777900528
$ws.Range("F9").Value = 4.262716023459958
$ws.Range("G9").Value = 19

$ws.Range("B10").Value = -0.5965755513627876
$ws.Range("C10").Value = 2.793546223734767
$ws.Range("D10").Value = 23.31997394578584
$ws.Range("E10").Value = 4.829075889420857
$ws.Range("F10").Value = 5.005168770995406
$ws.Range("G10").Value = 12

$ws.Range("B11").Value = -3.259825893570965
$ws.Range("C11").Value = 3.570597827222227
$ws.Range("D11").Value = 29.46234283587567
$ws.Range("E11").Value = 5.427922515647738
$ws.Range("F11").Value = 4.852303316400357
